# Update res_bus/vm_pu.xlsx results for the 380 kV case:
# slack-bus voltage setpoint (column B) changes from 1.05 to 1.02 p.u.,
# and the dependent bus voltage-magnitude results (columns C,D,E,I,J,K,L,N)
# are refreshed with the corresponding re-solved load-flow values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027028756080297
$ws.Range("D2").Value = 1.031212491392802
$ws.Range("E2").Value = 1.027184337646082
$ws.Range("I2").Value = 1.034222890360825
$ws.Range("J2").Value = 1.032188907551269
$ws.Range("K2").Value = 1.034021261770401
$ws.Range("L2").Value = 1.03000480634619
$ws.Range("N2").Value = 1.033654734193721

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027803685378983
$ws.Range("D3").Value = 1.031773704430455
$ws.Range("E3").Value = 1.02783739651774
$ws.Range("I3").Value = 1.034386620114168
$ws.Range("J3").Value = 1.032604745615784
$ws.Range("K3").Value = 1.034391655004504
$ws.Range("L3").Value = 1.030465955215684
$ws.Range("N3").Value = 1.034071162795978

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028305508122105
$ws.Range("D4").Value = 1.032137061397323
$ws.Range("E4").Value = 1.028260708966881
$ws.Range("I4").Value = 1.034491327511651
$ws.Range("J4").Value = 1.032873536411004
$ws.Range("K4").Value = 1.034630834484209
$ws.Range("L4").Value = 1.030764404393195
$ws.Range("N4").Value = 1.034340335304958

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028516566294637
$ws.Range("D5").Value = 1.032289866134371
$ws.Range("E5").Value = 1.028438844636566
$ws.Range("I5").Value = 1.034535049542256
$ws.Range("J5").Value = 1.032986466635506
$ws.Range("K5").Value = 1.034731266931604
$ws.Range("L5").Value = 1.0308898841184
$ws.Range("N5").Value = 1.034453425903331

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028552009224531
$ws.Range("D6").Value = 1.032315525567332
$ws.Range("E6").Value = 1.028468764610936
$ws.Range("I6").Value = 1.034542373223841
$ws.Range("J6").Value = 1.033005424005825
$ws.Range("K6").Value = 1.0347481229865
$ws.Range("L6").Value = 1.030910953358824
$ws.Range("N6").Value = 1.034472410195292

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02830832793302
$ws.Range("D7").Value = 1.032139102989732
$ws.Range("E7").Value = 1.028263088538813
$ws.Range("I7").Value = 1.034491912895528
$ws.Range("J7").Value = 1.032875045664722
$ws.Range("K7").Value = 1.03463217693525
$ws.Range("L7").Value = 1.030766081015354
$ws.Range("N7").Value = 1.034341846701989

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027290564840269
$ws.Range("D8").Value = 1.031402110385325
$ws.Range("E8").Value = 1.027404887610115
$ws.Range("I8").Value = 1.034278479129495
$ws.Range("J8").Value = 1.032329499997596
$ws.Range("K8").Value = 1.034146538473983
$ws.Range("L8").Value = 1.030160641483019
$ws.Range("N8").Value = 1.033795526297441

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025500218810762
$ws.Range("D9").Value = 1.030105163073955
$ws.Range("E9").Value = 1.025898379919762
$ws.Range("I9").Value = 1.033892949238903
$ws.Range("J9").Value = 1.031366062748962
$ws.Range("K9").Value = 1.033287091245231
$ws.Range("L9").Value = 1.029094266481848
$ws.Range("N9").Value = 1.032830720857406

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024308838945009
$ws.Range("D10").Value = 1.029241806907891
$ws.Range("E10").Value = 1.024898026078078
$ws.Range("I10").Value = 1.033629640545774
$ws.Range("J10").Value = 1.030722430560646
$ws.Range("K10").Value = 1.03271172676893
$ws.Range("L10").Value = 1.028383765366788
$ws.Range("N10").Value = 1.032186174637558

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023793499869064
$ws.Range("D11").Value = 1.02886829069914
$ws.Range("E11").Value = 1.024465831743695
$ws.Range("I11").Value = 1.033514144518104
$ws.Range("J11").Value = 1.030443430030799
$ws.Range("K11").Value = 1.032462035613863
$ws.Range("L11").Value = 1.028076227746639
$ws.Range("N11").Value = 1.031906777894957

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023602162450337
$ws.Range("D12").Value = 1.028729600581146
$ws.Range("E12").Value = 1.024305442596741
$ws.Range("I12").Value = 1.033471022435472
$ws.Range("J12").Value = 1.030339752599397
$ws.Range("K12").Value = 1.032369207268754
$ws.Range("L12").Value = 1.027962013223487
$ws.Range("N12").Value = 1.031802953229709

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023643201214335
$ws.Range("D13").Value = 1.028759347768825
$ws.Range("E13").Value = 1.024339839925513
$ws.Range("I13").Value = 1.033480282287515
$ws.Range("J13").Value = 1.030361993738743
$ws.Range("K13").Value = 1.03238912294938
$ws.Range("L13").Value = 1.027986511758658
$ws.Range("N13").Value = 1.031825225954024

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023777682165259
$ws.Range("D14").Value = 1.028856825490658
$ws.Range("E14").Value = 1.024452570914884
$ws.Range("I14").Value = 1.033510584551131
$ws.Range("J14").Value = 1.03043486091585
$ws.Range("K14").Value = 1.032454364057326
$ws.Range("L14").Value = 1.02806678634953
$ws.Range("N14").Value = 1.031898196610882

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023860551300809
$ws.Range("D15").Value = 1.028916891473283
$ws.Range("E15").Value = 1.024522047746853
$ws.Range("I15").Value = 1.033529225422765
$ws.Range("J15").Value = 1.030479750971645
$ws.Range("K15").Value = 1.032494550441538
$ws.Range("L15").Value = 1.02811624868964
$ws.Range("N15").Value = 1.031943150415704

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024343051762257
$ws.Range("D16").Value = 1.029266602914157
$ws.Range("E16").Value = 1.024926729906745
$ws.Range("I16").Value = 1.033637274495902
$ws.Range("J16").Value = 1.03074094063049
$ws.Range("K16").Value = 1.032728286382964
$ws.Range("L16").Value = 1.028404178157446
$ws.Range("N16").Value = 1.032204710993824

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024645856463726
$ws.Range("D17").Value = 1.029486055438039
$ws.Range("E17").Value = 1.025180836248098
$ws.Range("I17").Value = 1.033704654774372
$ws.Range("J17").Value = 1.030904697718892
$ws.Range("K17").Value = 1.032874755273379
$ws.Range("L17").Value = 1.028584820305422
$ws.Range("N17").Value = 1.032368700636075

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024822528938316
$ws.Range("D18").Value = 1.029614089388658
$ws.Range("E18").Value = 1.025329145197571
$ws.Range("I18").Value = 1.03374381353628
$ws.Range("J18").Value = 1.031000185077113
$ws.Range("K18").Value = 1.032960134501608
$ws.Range("L18").Value = 1.028690196717526
$ws.Range("N18").Value = 1.032464323597297

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024882778383819
$ws.Range("D19").Value = 1.029657750856931
$ws.Range("E19").Value = 1.025379730432923
$ws.Range("I19").Value = 1.033757141377546
$ws.Range("J19").Value = 1.03103273877275
$ws.Range("K19").Value = 1.03298923747892
$ws.Range("L19").Value = 1.028726129154387
$ws.Range("N19").Value = 1.032496923522916

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024613363032513
$ws.Range("D20").Value = 1.02946250703791
$ws.Range("E20").Value = 1.025153563406471
$ws.Range("I20").Value = 1.033697440294471
$ws.Range("J20").Value = 1.030887131162044
$ws.Range("K20").Value = 1.032859046081479
$ws.Range("L20").Value = 1.028565437973856
$ws.Range("N20").Value = 1.032351109132701

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023738078585587
$ws.Range("D21").Value = 1.028828119323329
$ws.Range("E21").Value = 1.024419370372151
$ws.Range("I21").Value = 1.033501667403112
$ws.Range("J21").Value = 1.030413404557663
$ws.Range("K21").Value = 1.032435154427795
$ws.Range("L21").Value = 1.02804314695661
$ws.Range("N21").Value = 1.031876709782204

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023188230535262
$ws.Range("D22").Value = 1.028429547337334
$ws.Range("E22").Value = 1.023958605982246
$ws.Range("I22").Value = 1.033377294874859
$ws.Range("J22").Value = 1.030115298728812
$ws.Range("K22").Value = 1.032168163565304
$ws.Range("L22").Value = 1.027714870465271
$ws.Range("N22").Value = 1.031578180608883

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023479669447885
$ws.Range("D23").Value = 1.028640809489753
$ws.Range("E23").Value = 1.024202784478957
$ws.Range("I23").Value = 1.03344334836204
$ws.Range("J23").Value = 1.030273353986425
$ws.Range("K23").Value = 1.032309744881542
$ws.Range("L23").Value = 1.027888885248066
$ws.Range("N23").Value = 1.031736460323092

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024628045256168
$ws.Range("D24").Value = 1.029473147450751
$ws.Range("E24").Value = 1.025165886542509
$ws.Range("I24").Value = 1.033700700649394
$ws.Range("J24").Value = 1.0308950688239
$ws.Range("K24").Value = 1.032866144554888
$ws.Range("L24").Value = 1.028574195981694
$ws.Range("N24").Value = 1.032359058066947

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025962689370182
$ws.Range("D25").Value = 1.030440238550345
$ws.Range("E25").Value = 1.026287154521877
$ws.Range("I25").Value = 1.033993730272588
$ws.Range("J25").Value = 1.031615376673928
$ws.Range("K25").Value = 1.03350970786159
$ws.Range("L25").Value = 1.029369883185305
$ws.Range("N25").Value = 1.033080388836742
